$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet: "Property1" -> "DataNode"
$ws.Name = "DataNode"

# Update the active selection (cosmetic UI state) from A9 to E23
[void]$ws.Range("E23").Select()

# Nudge column widths for columns A and E very slightly (manual border-drag edit)
$ws.Columns.Item(1).ColumnWidth = 11.785714285714285
$ws.Columns.Item(5).ColumnWidth = 13.071428571428571
